# Fruta / hortaliza, semanal
# Insert a new weekly record at row 32 (pushing the existing rows 32-38
# down to 33-39) for "Vega Monumental Concepción - Poroto granado".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 32:38 down to 33:39, leaving a blank row 32 to populate.
$ws.Rows.Item(32).Insert()

# Populate the new row 32 with the latest weekly price entry.
$ws.Range("A32").Value = 11
$ws.Range("B32").Value = "Vega Monumental Concepción"
$ws.Range("C32").Value = "Bíobío"
$ws.Range("D32").Value = 44642
$ws.Range("E32").Value = 8
$ws.Range("F32").Value = 100112030
$ws.Range("G32").Value = "Poroto granado"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 120
$ws.Range("K32").Value = 20000
$ws.Range("L32").Value = 20000
$ws.Range("M32").Value = 20000
$ws.Range("N32").Value = "$/saco 25 kilos"
$ws.Range("O32").Value = "Región Metropolitana"
$ws.Range("P32").Value = 800
$ws.Range("Q32").Value = 25
$ws.Range("R32").Value = "Hortaliza"
